$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "SFX_1" entry (id 2000) with "None"
$ws.Range("B2").Value = "None"

# Make sure the main scenes / active selection is set to C5 (single cell)
[void]$ws.Range("C5").Select()
